$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A45").Value = 44
$ws.Range("B45").Value = "TheRiverMan"
$ws.Range("C45").Value = "The River Man"
$ws.Range("D45").Value = 2
$ws.Range("E45").Value = "Purgatory"
$ws.Range("I45").Value = 1872
$ws.Range("J45").Value = 816
